$wb = $excel.ActiveWorkbook

# --- Sheet: Top 50 Cryptocurrencies ---
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$ws1.Range("C2").Value = 98973
$ws1.Range("D2").Value = 1960031497170
$ws1.Range("E2").Value = 112063105617
$ws1.Range("F2").Value = 1.34469
$ws1.Range("C3").Value = 3391.73
$ws1.Range("D3").Value = 408586496268
$ws1.Range("E3").Value = 57188238984
$ws1.Range("F3").Value = 8.5235
$ws1.Range("D4").Value = 130911146723
$ws1.Range("E4").Value = 161208017217
$ws1.Range("F4").Value = -0.06610000000000001
$ws1.Range("C5").Value = 261.97
$ws1.Range("D5").Value = 124520979419
$ws1.Range("E5").Value = 14956566986
$ws1.Range("F5").Value = 8.4557
$ws1.Range("C6").Value = 634.28
$ws1.Range("D6").Value = 92589848319
$ws1.Range("E6").Value = 2461265568
$ws1.Range("F6").Value = 3.64408
$ws1.Range("C7").Value = 1.4
$ws1.Range("D7").Value = 79352270388
$ws1.Range("E7").Value = 18051002346
$ws1.Range("F7").Value = 24.9139
$ws1.Range("C8").Value = 0.396461
$ws1.Range("D8").Value = 58243325207
$ws1.Range("E8").Value = 9794203793
$ws1.Range("F8").Value = 2.5468
$ws1.Range("D9").Value = 38324646891
$ws1.Range("E9").Value = 13566669467
$ws1.Range("F9").Value = -0.03895
$ws1.Range("C10").Value = 3393.11
$ws1.Range("D10").Value = 33237192980
$ws1.Range("E10").Value = 147411661
$ws1.Range("F10").Value = 8.583399999999999
$ws1.Range("C11").Value = 0.883182
$ws1.Range("D11").Value = 31606335117
$ws1.Range("E11").Value = 3575419706
$ws1.Range("F11").Value = 11.9422
$ws1.Range("C12").Value = 0.200351
$ws1.Range("D12").Value = 17309848780
$ws1.Range("E12").Value = 1067883664
$ws1.Range("F12").Value = 1.44964
$ws1.Range("C13").Value = 36.41
$ws1.Range("D13").Value = 14894009679
$ws1.Range("E13").Value = 1051408833
$ws1.Range("F13").Value = 6.92825
$ws1.Range("C14").Value = 0.00002501
$ws1.Range("D14").Value = 14744333428
$ws1.Range("E14").Value = 1608726152
$ws1.Range("F14").Value = 3.15487
$ws1.Range("C15").Value = 4008.33
$ws1.Range("D15").Value = 14474546253
$ws1.Range("E15").Value = 167732168
$ws1.Range("F15").Value = 8.51431
$ws1.Range("C16").Value = 98469
$ws1.Range("D16").Value = 14383538226
$ws1.Range("E16").Value = 845619925
$ws1.Range("F16").Value = 1.4294
$ws1.Range("D17").Value = 14174777304
$ws1.Range("E17").Value = 637789593
$ws1.Range("F17").Value = 3.41937
$ws1.Range("C18").Value = 3.61
$ws1.Range("D18").Value = 10273657188
$ws1.Range("E18").Value = 2206148805
$ws1.Range("F18").Value = 0.14521
$ws1.Range("C19").Value = 497.52
$ws1.Range("D19").Value = 9852940247
$ws1.Range("E19").Value = 1909542516
$ws1.Range("F19").Value = -4.17898
$ws1.Range("C20").Value = 3391.08
$ws1.Range("D20").Value = 9719868552
$ws1.Range("E20").Value = 1466024154
$ws1.Range("F20").Value = 8.617419999999999
$ws1.Range("C21").Value = 15.3
$ws1.Range("D21").Value = 9590618652
$ws1.Range("E21").Value = 1256306565
$ws1.Range("F21").Value = 4.37437
$ws1.Range("A22").Value = "Pepe"
$ws1.Range("B22").Value = "pepe"
$ws1.Range("C22").Value = 0.00002138
$ws1.Range("D22").Value = 8992937206
$ws1.Range("E22").Value = 6818526335
$ws1.Range("F22").Value = 9.91516
$ws1.Range("A23").Value = "Polkadot"
$ws1.Range("B23").Value = "dot"
$ws1.Range("C23").Value = 6.23
$ws1.Range("D23").Value = 8988502483
$ws1.Range("E23").Value = 830212815
$ws1.Range("F23").Value = 9.276669999999999
$ws1.Range("C24").Value = 0.283394
$ws1.Range("D24").Value = 8500667311
$ws1.Range("E24").Value = 2299576081
$ws1.Range("F24").Value = 19.34584
$ws1.Range("C25").Value = 8.789999999999999
$ws1.Range("D25").Value = 8129859866
$ws1.Range("E25").Value = 3437669
$ws1.Range("F25").Value = 3.23414
$ws1.Range("C26").Value = 5.81
$ws1.Range("D26").Value = 7079430306
$ws1.Range("E26").Value = 1013157831
$ws1.Range("F26").Value = 4.4035
$ws1.Range("C27").Value = 90.98999999999999
$ws1.Range("D27").Value = 6843655404
$ws1.Range("E27").Value = 1421163417
$ws1.Range("F27").Value = 4.39455
$ws1.Range("C28").Value = 12.15
$ws1.Range("D28").Value = 6471704530
$ws1.Range("E28").Value = 863170407
$ws1.Range("F28").Value = 3.58351
$ws1.Range("C29").Value = 3572.68
$ws1.Range("D29").Value = 6215050475
$ws1.Range("E29").Value = 105919148
$ws1.Range("F29").Value = 8.699859999999999
$ws1.Range("C30").Value = 9.449999999999999
$ws1.Range("D30").Value = 5675174761
$ws1.Range("E30").Value = 867407640
$ws1.Range("F30").Value = 6.44424
$ws1.Range("C31").Value = 0.203407
$ws1.Range("D31").Value = 5518803426
$ws1.Range("E31").Value = 130186401
$ws1.Range("F31").Value = 15.95664
$ws1.Range("C32").Value = 0.99739
$ws1.Range("D32").Value = 5219437344
$ws1.Range("E32").Value = 91138
$ws1.Range("F32").Value = -0.41544
$ws1.Range("C33").Value = 0.132838
$ws1.Range("D33").Value = 5102360242
$ws1.Range("E33").Value = 902012626
$ws1.Range("F33").Value = 5.76491
$ws1.Range("C34").Value = 9.65
$ws1.Range("D34").Value = 4574959826
$ws1.Range("E34").Value = 274551581
$ws1.Range("F34").Value = 6.2297
$ws1.Range("C35").Value = 28.08
$ws1.Range("D35").Value = 4200949688
$ws1.Range("E35").Value = 877961705
$ws1.Range("F35").Value = 5.00049
$ws1.Range("C36").Value = 0.0000521
$ws1.Range("D36").Value = 3909767343
$ws1.Range("E36").Value = 1680780764
$ws1.Range("F36").Value = 2.43461
$ws1.Range("A37").Value = "Render"
$ws1.Range("B37").Value = "render"
$ws1.Range("C37").Value = 7.4
$ws1.Range("D37").Value = 3830524194
$ws1.Range("E37").Value = 435811535
$ws1.Range("F37").Value = 0.63923
$ws1.Range("A38").Value = "Kaspa"
$ws1.Range("B38").Value = "kas"
$ws1.Range("C38").Value = 0.151489
$ws1.Range("D38").Value = 3821068085
$ws1.Range("E38").Value = 151494380
$ws1.Range("F38").Value = -0.29508
$ws1.Range("C39").Value = 0.473954
$ws1.Range("D39").Value = 3779000218
$ws1.Range("E39").Value = 492239332
$ws1.Range("F39").Value = 7.932
$ws1.Range("C40").Value = 507.23
$ws1.Range("D40").Value = 3743956169
$ws1.Range("E40").Value = 283850524
$ws1.Range("F40").Value = 2.43937
$ws1.Range("D41").Value = 3689800135
$ws1.Range("E41").Value = 224694439
$ws1.Range("F41").Value = -0.05866
$ws1.Range("C42").Value = 24.81
$ws1.Range("D42").Value = 3577004631
$ws1.Range("E42").Value = 33482480
$ws1.Range("F42").Value = 2.55998
$ws1.Range("C43").Value = 0.999814
$ws1.Range("D43").Value = 3444066515
$ws1.Range("E43").Value = 154806556
$ws1.Range("F43").Value = -0.09016
$ws1.Range("C44").Value = 3.4
$ws1.Range("D44").Value = 3387805331
$ws1.Range("E44").Value = 1284980672
$ws1.Range("F44").Value = 6.28501
$ws1.Range("C45").Value = 3.72
$ws1.Range("D45").Value = 3351839645
$ws1.Range("E45").Value = 301332116
$ws1.Range("F45").Value = 3.21506
$ws1.Range("D46").Value = 3344641105
$ws1.Range("E46").Value = 485928047
$ws1.Range("F46").Value = 2.65695
$ws1.Range("C47").Value = 0.792956
$ws1.Range("D47").Value = 3248313579
$ws1.Range("E47").Value = 1672915093
$ws1.Range("F47").Value = 13.78308
$ws1.Range("C48").Value = 161.3
$ws1.Range("D48").Value = 2977336240
$ws1.Range("E48").Value = 84382065
$ws1.Range("F48").Value = -0.91428
$ws1.Range("D49").Value = 2948469908
$ws1.Range("E49").Value = 364748635
$ws1.Range("F49").Value = 1.20934
$ws1.Range("C50").Value = 0.847034
$ws1.Range("D50").Value = 2856006968
$ws1.Range("E50").Value = 185488113
$ws1.Range("F50").Value = 15.42842
$ws1.Range("C51").Value = 4.72
$ws1.Range("D51").Value = 2833166913
$ws1.Range("E51").Value = 579787658
$ws1.Range("F51").Value = 7.26919

# --- Sheet: Top 5 by Market Cap ---
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")

$ws2.Range("B2").Value = 1960031497170
$ws2.Range("B3").Value = 408586496268
$ws2.Range("B4").Value = 130911146723
$ws2.Range("B5").Value = 124520979419
$ws2.Range("B6").Value = 92589848319

# --- Sheet: Summary ---
$ws3 = $wb.Worksheets.Item("Summary")

$ws3.Range("B2").Value = "`$4351.03"
$ws3.Range("B3").Value = "XRP (24.91%)"
$ws3.Range("B4").Value = "Bitcoin Cash (-4.18%)"
